# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1740.4
$ws.Range("I18").Value = 1300
$ws.Range("J18").Value = 3502
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 3502
$ws.Range("M18").Value = -1016
$ws.Range("N18").Value = -4070
$ws.Range("H20").Value = 3344833
$ws.Range("I20").Value = 3344833
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3344833
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -3344603
$ws.Range("H33").Value = 229.75
$ws.Range("I33").Value = 246.93333
$ws.Range("K33").Value = 246.93333
$ws.Range("M33").Value = -17.93333000000001
$ws.Range("H35").Value = 3344833
$ws.Range("I35").Value = 3344833
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3344833
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -3344454
$ws.Range("H40").Value = 27165.5
$ws.Range("J40").Value = 16760.572
$ws.Range("L40").Value = 16760.572
$ws.Range("N40").Value = -17110.572
$ws.Range("H64").Value = 5549.1665
$ws.Range("J64").Value = 5559
$ws.Range("L64").Value = 5559
$ws.Range("N64").Value = -6055
$ws.Range("H67").Value = 5549.1665
$ws.Range("J67").Value = 5559
$ws.Range("L67").Value = 5559
$ws.Range("N67").Value = -7275
$ws.Range("H69").Value = 55562140
$ws.Range("J69").Value = 7899.8
$ws.Range("L69").Value = 23699.4
$ws.Range("N69").Value = -25447.4
$ws.Range("H72").Value = 55562140
$ws.Range("J72").Value = 7899.8
$ws.Range("L72").Value = 71098.2
$ws.Range("N72").Value = -79834.2
$ws.Range("H94").Value = 7282.385
$ws.Range("I94").Value = 3769.4443
$ws.Range("K94").Value = 3769.4443
$ws.Range("M94").Value = -3318.4443
$ws.Range("H98").Value = 3390.0435
$ws.Range("I98").Value = 3357.5264
$ws.Range("K98").Value = 3357.5264
$ws.Range("M98").Value = -1859.5264
$ws.Range("H107").Value = 881.8570999999999
$ws.Range("I107").Value = 840.44446
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 840.44446
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1079.55554
$ws.Range("N107").Value = -5840
$ws.Range("H112").Value = 2841.8462
$ws.Range("J112").Value = 2862.6667
$ws.Range("L112").Value = 8588.000100000001
$ws.Range("N112").Value = -10804.0001
$ws.Range("H116").Value = 1281756.5
$ws.Range("I116").Value = 2754753.5
$ws.Range("K116").Value = 2754753.5
$ws.Range("M116").Value = -2751311.5
$ws.Range("H122").Value = 3390.0435
$ws.Range("I122").Value = 3357.5264
$ws.Range("K122").Value = 10072.5792
$ws.Range("M122").Value = -7622.5792
$ws.Range("H125").Value = 5888.5625
$ws.Range("I125").Value = 2666.6667
$ws.Range("K125").Value = 24000.0003
$ws.Range("M125").Value = -21540.0003
$ws.Range("H135").Value = 4390.3335
$ws.Range("I135").Value = 1211.25
$ws.Range("J135").Value = 10748.5
$ws.Range("K135").Value = 10901.25
$ws.Range("L135").Value = 96736.5
$ws.Range("M135").Value = -8366.25
$ws.Range("N135").Value = -101806.5
$ws.Range("H137").Value = 9037598
$ws.Range("I137").Value = 557109.5600000001
$ws.Range("K137").Value = 1671328.68
$ws.Range("M137").Value = -1668778.68
$ws.Range("H139").Value = 134606.62
$ws.Range("I139").Value = 131000
$ws.Range("J139").Value = 134847.06
$ws.Range("K139").Value = 131000
$ws.Range("L139").Value = 134847.06
$ws.Range("M139").Value = -125860
$ws.Range("N139").Value = -145127.06

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13160382
$ws.Range("I74").Value = 15627268
$ws.Range("J74").Value = 3658
$ws.Range("K74").Value = 15627268
$ws.Range("L74").Value = 3658
$ws.Range("M74").Value = -15626394
$ws.Range("N74").Value = -5406
$ws.Range("H77").Value = 13160382
$ws.Range("I77").Value = 15627268
$ws.Range("J77").Value = 3658
$ws.Range("K77").Value = 78136340
$ws.Range("L77").Value = 18290
$ws.Range("M77").Value = -78131972
$ws.Range("N77").Value = -27026

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 90000
$ws.Range("J57").Value = 90000
$ws.Range("L57").Value = 90000
$ws.Range("N57").Value = -91440
$ws.Range("H132").Value = 117011.12
$ws.Range("J132").Value = 117011.12
$ws.Range("L132").Value = 117011.12
$ws.Range("N132").Value = -127131.12
$ws.Range("H136").Value = 90000
$ws.Range("J136").Value = 90000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -100200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19610420
$ws.Range("I31").Value = 25001776
$ws.Range("K31").Value = 25001776
$ws.Range("M31").Value = -25001481
$ws.Range("H34").Value = 19610420
$ws.Range("I34").Value = 25001776
$ws.Range("K34").Value = 25001776
$ws.Range("M34").Value = -25001574
$ws.Range("H62").Value = 35385.85
$ws.Range("I62").Value = 32293.467
$ws.Range("J62").Value = 39251.332
$ws.Range("K62").Value = 32293.467
$ws.Range("L62").Value = 39251.332
$ws.Range("M62").Value = -31669.467
$ws.Range("N62").Value = -40499.332
$ws.Range("H65").Value = 35385.85
$ws.Range("I65").Value = 32293.467
$ws.Range("J65").Value = 39251.332
$ws.Range("K65").Value = 161467.335
$ws.Range("L65").Value = 196256.66
$ws.Range("M65").Value = -158347.335
$ws.Range("N65").Value = -202496.66
$ws.Range("H132").Value = 43480600
$ws.Range("I132").Value = 52633800
$ws.Range("J132").Value = 2902
$ws.Range("K132").Value = 157901400
$ws.Range("L132").Value = 8706
$ws.Range("M132").Value = -157898870
$ws.Range("N132").Value = -13766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2112
$ws.Range("I92").Value = 2222
$ws.Range("K92").Value = 6666
$ws.Range("M92").Value = -5418
$ws.Range("H132").Value = 2958
$ws.Range("I132").Value = 2050
$ws.Range("K132").Value = 18450
$ws.Range("M132").Value = -15920
$ws.Range("H134").Value = 8083.3
$ws.Range("I134").Value = 5729.125
$ws.Range("K134").Value = 17187.375
$ws.Range("M134").Value = -12117.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23615.385
$ws.Range("J46").Value = 25000
$ws.Range("L46").Value = 25000
$ws.Range("N46").Value = -25312
$ws.Range("H102").Value = 23186926
$ws.Range("I102").Value = 31879714
$ws.Range("K102").Value = 31879714
$ws.Range("M102").Value = -31878092
$ws.Range("H122").Value = 3454.8936
$ws.Range("I122").Value = 1715.1818
$ws.Range("J122").Value = 4985.84
$ws.Range("K122").Value = 5145.5454
$ws.Range("L122").Value = 14957.52
$ws.Range("M122").Value = -2695.5454
$ws.Range("N122").Value = -19857.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6845
$ws.Range("I7").Value = 4750
$ws.Range("J7").Value = 7310.5557
$ws.Range("K7").Value = 4750
$ws.Range("L7").Value = 7310.5557
$ws.Range("M7").Value = -4638
$ws.Range("N7").Value = -7534.5557
$ws.Range("H40").Value = 29173556
$ws.Range("I40").Value = 20838396
$ws.Range("K40").Value = 20838396
$ws.Range("M40").Value = -20838260
$ws.Range("H82").Value = 3895.7036
$ws.Range("I82").Value = 3021.0625
$ws.Range("K82").Value = 3021.0625
$ws.Range("M82").Value = -2660.0625
$ws.Range("H85").Value = 3895.7036
$ws.Range("I85").Value = 3021.0625
$ws.Range("K85").Value = 3021.0625
$ws.Range("M85").Value = -1773.0625
$ws.Range("H93").Value = 2712.4666
$ws.Range("I93").Value = 1927.6666
$ws.Range("K93").Value = 1927.6666
$ws.Range("M93").Value = -679.6666
$ws.Range("H108").Value = 105999.4
$ws.Range("J108").Value = 105999.4
$ws.Range("L108").Value = 105999.4
$ws.Range("N108").Value = -113679.4
$ws.Range("H126").Value = 6845
$ws.Range("I126").Value = 4750
$ws.Range("J126").Value = 7310.5557
$ws.Range("K126").Value = 14250
$ws.Range("L126").Value = 21931.6671
$ws.Range("M126").Value = -11780
$ws.Range("N126").Value = -26871.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2467.5652
$ws.Range("I122").Value = 2308.889
$ws.Range("J122").Value = 3038.8
$ws.Range("K122").Value = 6926.667
$ws.Range("L122").Value = 9116.400000000001
$ws.Range("M122").Value = -4476.667
$ws.Range("N122").Value = -14016.4
$ws.Range("H126").Value = 333334500
$ws.Range("I126").Value = 333334500
$ws.Range("K126").Value = 1000003500
$ws.Range("M126").Value = -1000001030
